# Update the dSF (column F) values for a set of re-pulled rows.
# Mapping: row number (1-based) -> new value
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    5  = -2
    8  = -3
    14 = -2
    16 = -3
    17 = -2
    21 = -4
    22 = 1
    23 = -3
    29 = 4
    30 = -10
    36 = 4
    40 = -4
    41 = 3
    42 = -5
    44 = -1
    50 = 8
    52 = 0
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
